# Remove the "Hi UI.Name. " greeting that used the UI.Name merge-field
# style variable, leaving just a single space in its place (matching the
# author's change of getting rid of the UI.variable placeholders).
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("Hi UI.Name. ", $true, $false, $false, $false, $false, `
               $true, 1, $false, " ", 2)
